$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (border/bold/alignment) of the existing year-label
# cell A6 onto the two new year-label cells A7 and A8 so they keep the
# same look as A2:A6.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

# Row 7: 2021年 - full set of values
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 98.7
$ws.Range("C7").Value = 100.8
$ws.Range("D7").Value = 100
$ws.Range("E7").Value = 102.2
$ws.Range("F7").Value = 101.1
$ws.Range("G7").Value = 100.2
$ws.Range("H7").Value = 100.4

# Row 8: 2022年 - only the H column (生活用品及服务类) figure is available so far
$ws.Range("A8").Value = "2022年"
$ws.Range("H8").Value = 101
